$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column G
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-13 20:27:34"

    # zh-cn sheet: "Latest Handoff Datetime" column H, and "Priority" column E
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 20:27:26"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"

    # de-de sheet: "Priority" column E
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
